# seminario: atualizando slide conclusao
#
# 1) Refresh the cached "datetimeFigureOut" placeholder text (03/06/2024 -> 04/06/2024)
#    on every slide layout, the slide master and the notes master.
# 2) Rewrite the three conclusion bullets on the "Conclusoes" slide.
# 3) Tidy up the page-number shape on the last slide (merge the stray
#    endParaRPr run into the visible run).

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($shapes, $newText)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

# -- 1. Date placeholders -----------------------------------------------
Update-DatePlaceholder $p.SlideMaster.Shapes "04/06/2024"

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes "04/06/2024"
}

Update-DatePlaceholder $p.NotesMaster.Shapes "04/06/2024"

# -- 2. Conclusion slide bullets -----------------------------------------
$slide = $p.Slides.Item(26)
$body = $slide.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

# Paragraph 1: split into 3 runs, middle one flagged as a foreign word.
$para1 = $tr.Paragraphs(1)
$para1.Text = "#"
$para1 = $tr.Paragraphs(1)
$para1.Text = "O algoritmo AHC com single "
$run1b = $tr.Paragraphs(1).InsertAfter("linkage")
$run1c = $run1b.InsertAfter(" demonstrou resultados de agrupamento est" + [char]0xE1 + "veis")

# Paragraph 2: append extra clause, with a foreign word run in the middle.
$para2 = $tr.Paragraphs(2)
$run2a = $para2.InsertAfter(" e s" + [char]0xE3 + "o adequados para determinar o n" + [char]0xFA + "mero ideal de agrupamentos para dados lineares, ")
$run2b = $run2a.InsertAfter("manifold")
$run2c = $run2b.InsertAfter(", anulares e convexos.")

# Paragraph 3: simple word swap (otimizacao -> adaptacao).
$para3 = $tr.Paragraphs(3)
$para3.Text = "#"
$para3 = $tr.Paragraphs(3)
$para3.Text = "O algoritmo ONCD " + [char]0xE9 + " extens" + [char]0xED + "vel " + [char]0xE0 + " outros algoritmos de agrupamento (sendo necess" + [char]0xE1 + "rio certo grau de adapta" + [char]0xE7 + [char]0xE3 + "o)"

# -- 3. Page-number shape on the closing slide ---------------------------
$lastSlide = $p.Slides.Item(27)
$pageShape = $lastSlide.Shapes.Item(5)
$pageRange = $pageShape.TextFrame.TextRange
$pageRange.Paragraphs(1).Text = "27"
$pageRange.Font.Size = 1200
$pageRange.Font.Color = $pageRange.Font.Color
